# Generate Report for Handoff
#
# The handoff/report generator refreshed the "Ready for handoff" timestamps
# for the 564bd0cf-42e9-4340-a0e9-fb94fd5c91e4 file across the Overview
# sheet (Latest HO Xliff Generate Date) and each language sheet's
# "Latest Handoff Datetime" column (row for that file is row 4 on every
# sheet).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: column G is "Latest HO Xliff Generate Date" ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4").Value = "2016-10-27 09:07:21"

# --- zh-cn sheet: column H is "Latest Handoff Datetime" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H4").Value = "2016-10-27 09:07:10"

# --- de-de sheet: column H is "Latest Handoff Datetime" ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H4").Value = "2016-10-27 09:07:21"
